$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 937.4
$ws.Range("I2").Value = 471.75
$ws.Range("J2").Value = 2800
$ws.Range("K2").Value = 471.75
$ws.Range("L2").Value = 2800
$ws.Range("M2").Value = -358.75
$ws.Range("N2").Value = -3026
$ws.Range("H12").Value = 460.53845
$ws.Range("I12").Value = 124.75
$ws.Range("J12").Value = 997.8
$ws.Range("K12").Value = 124.75
$ws.Range("L12").Value = 997.8
$ws.Range("M12").Value = 45.25
$ws.Range("N12").Value = -1337.8
$ws.Range("H51").Value = 5749.2
$ws.Range("J51").Value = 6249.25
$ws.Range("L51").Value = 6249.25
$ws.Range("N51").Value = -7217.25
$ws.Range("H61").Value = 988.1429
$ws.Range("I61").Value = 1086.1666
$ws.Range("J61").Value = 400
$ws.Range("K61").Value = 3258.4998
$ws.Range("L61").Value = 1200
$ws.Range("M61").Value = -3086.4998
$ws.Range("N61").Value = -1544
$ws.Range("H99").Value = 2234.818
$ws.Range("I99").Value = 1097.8
$ws.Range("J99").Value = 3182.3333
$ws.Range("K99").Value = 3293.4
$ws.Range("L99").Value = 9546.999899999999
$ws.Range("M99").Value = -1795.4
$ws.Range("N99").Value = -12542.9999
$ws.Range("H101").Value = 7976.6
$ws.Range("I101").Value = 7029.6665
$ws.Range("J101").Value = 9397
$ws.Range("K101").Value = 21088.9995
$ws.Range("L101").Value = 28191
$ws.Range("M101").Value = -19466.9995
$ws.Range("N101").Value = -31435
$ws.Range("H103").Value = 415.6
$ws.Range("I103").Value = 314.5
$ws.Range("K103").Value = 943.5
$ws.Range("M103").Value = -357.5
$ws.Range("H127").Value = 730.25
$ws.Range("I127").Value = 526.3
$ws.Range("J127").Value = 1750
$ws.Range("K127").Value = 1578.9
$ws.Range("L127").Value = 5250
$ws.Range("M127").Value = 3381.1
$ws.Range("N127").Value = -15170
$ws.Range("H132").Value = 3078.8262
$ws.Range("I132").Value = 3226.4187
$ws.Range("J132").Value = 963.3333
$ws.Range("K132").Value = 9679.2561
$ws.Range("L132").Value = 2889.9999
$ws.Range("M132").Value = -7149.256100000001
$ws.Range("N132").Value = -7949.9999
$ws.Range("H138").Value = 3289.8416
$ws.Range("I138").Value = 3317.9033
$ws.Range("J138").Value = 3272.7844
$ws.Range("K138").Value = 9953.7099
$ws.Range("L138").Value = 9818.3532
$ws.Range("M138").Value = -4813.7099
$ws.Range("N138").Value = -20098.3532
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("H61").Value = 2457.077
$ws.Range("I61").Value = 2148.9092
$ws.Range("J61").Value = 4152
$ws.Range("K61").Value = 2148.9092
$ws.Range("L61").Value = 4152
$ws.Range("M61").Value = -1936.9092
$ws.Range("N61").Value = -4576
$ws.Range("H63").Value = 1257.6
$ws.Range("I63").Value = 1072.625
$ws.Range("J63").Value = 1997.5
$ws.Range("K63").Value = 1072.625
$ws.Range("L63").Value = 1997.5
$ws.Range("M63").Value = -386.625
$ws.Range("N63").Value = -3369.5
$ws.Range("H66").Value = 1257.6
$ws.Range("I66").Value = 1072.625
$ws.Range("J66").Value = 1997.5
$ws.Range("K66").Value = 5363.125
$ws.Range("L66").Value = 9987.5
$ws.Range("M66").Value = -1931.125
$ws.Range("N66").Value = -16851.5
$ws.Range("H74").Value = 1743
$ws.Range("I74").Value = 986.2
$ws.Range("J74").Value = 2499.8
$ws.Range("K74").Value = 986.2
$ws.Range("L74").Value = 2499.8
$ws.Range("M74").Value = -112.2
$ws.Range("N74").Value = -4247.8
$ws.Range("H77").Value = 1743
$ws.Range("I77").Value = 986.2
$ws.Range("J77").Value = 2499.8
$ws.Range("K77").Value = 4931
$ws.Range("L77").Value = 12499
$ws.Range("M77").Value = -563
$ws.Range("N77").Value = -21235
$ws.Range("H117").Value = 30000
$ws.Range("J117").Value = 30000
$ws.Range("L117").Value = 30000
$ws.Range("N117").Value = -39178
$ws.Range("H136").Value = 2457.077
$ws.Range("I136").Value = 2148.9092
$ws.Range("J136").Value = 4152
$ws.Range("K136").Value = 6446.7276
$ws.Range("L136").Value = 12456
$ws.Range("M136").Value = -3896.7276
$ws.Range("N136").Value = -17556
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 27999
$ws.Range("I35").Value = 27999
$ws.Range("K35").Value = 27999
$ws.Range("M35").Value = -27689
$ws.Range("H50").Value = 52233.332
$ws.Range("J50").Value = 73350
$ws.Range("L50").Value = 73350
$ws.Range("N50").Value = -74498
$ws.Range("H82").Value = 13922.385
$ws.Range("I82").Value = 3726.4546
$ws.Range("K82").Value = 3726.4546
$ws.Range("M82").Value = -3343.4546
$ws.Range("H85").Value = 13922.385
$ws.Range("I85").Value = 3726.4546
$ws.Range("K85").Value = 3726.4546
$ws.Range("M85").Value = -2400.4546
$ws.Range("H88").Value = 48333.332
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("H91").Value = 48333.332
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("H105").Value = 3388.3076
$ws.Range("I105").Value = 2280.375
$ws.Range("K105").Value = 2280.375
$ws.Range("M105").Value = -533.375
$ws.Range("H124").Value = 73665.664
$ws.Range("J124").Value = 73665.664
$ws.Range("L124").Value = 73665.664
$ws.Range("N124").Value = -83485.664
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2487.2144
$ws.Range("I58").Value = 2092.7
$ws.Range("J58").Value = 3473.5
$ws.Range("K58").Value = 2092.7
$ws.Range("L58").Value = 3473.5
$ws.Range("M58").Value = -1889.7
$ws.Range("N58").Value = -3879.5
$ws.Range("H103").Value = 26617.8
$ws.Range("I103").Value = 15772.5
$ws.Range("K103").Value = 15772.5
$ws.Range("M103").Value = -14600.5
$ws.Range("H136").Value = 2487.2144
$ws.Range("I136").Value = 2092.7
$ws.Range("J136").Value = 3473.5
$ws.Range("K136").Value = 6278.099999999999
$ws.Range("L136").Value = 10420.5
$ws.Range("M136").Value = -3728.099999999999
$ws.Range("N136").Value = -15520.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 799.7143
$ws.Range("I8").Value = 799.7143
$ws.Range("K8").Value = 2399.1429
$ws.Range("M8").Value = -2260.1429
$ws.Range("H74").Value = 10999.667
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 10999.667
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 32999.001
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -35121.001
$ws.Range("H77").Value = 10999.667
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 10999.667
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 98997.003
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -109605.003
$ws.Range("H81").Value = 12750
$ws.Range("I81").Value = 8000
$ws.Range("J81").Value = 17500
$ws.Range("K81").Value = 24000
$ws.Range("L81").Value = 52500
$ws.Range("M81").Value = -22877
$ws.Range("N81").Value = -54746
$ws.Range("H84").Value = 12750
$ws.Range("I84").Value = 8000
$ws.Range("J84").Value = 17500
$ws.Range("K84").Value = 72000
$ws.Range("L84").Value = 157500
$ws.Range("M84").Value = -66384
$ws.Range("N84").Value = -168732
$ws.Range("H101").Value = 25548.092
$ws.Range("I101").Value = 18250
$ws.Range("K101").Value = 54750
$ws.Range("M101").Value = -52316
$ws.Range("H131").Value = 21165356
$ws.Range("I131").Value = 10101968
$ws.Range("J131").Value = 33335082
$ws.Range("K131").Value = 30305904
$ws.Range("L131").Value = 100005246
$ws.Range("M131").Value = -30300864
$ws.Range("N131").Value = -100015326
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 6031.6665
$ws.Range("I113").Value = 5999
$ws.Range("J113").Value = 6038.2
$ws.Range("K113").Value = 5999
$ws.Range("L113").Value = 6038.2
$ws.Range("M113").Value = -3829
$ws.Range("N113").Value = -10378.2
$ws.Range("H122").Value = 2699.1667
$ws.Range("I122").Value = 2539
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 7617
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -5167
$ws.Range("N122").Value = -15400
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1829.0344
$ws.Range("I61").Value = 1779.3334
$ws.Range("K61").Value = 1779.3334
$ws.Range("M61").Value = -1577.3334
$ws.Range("H62").Value = 49800
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 49800
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H70").Value = 33499.4
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 33499.4
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 33499.4
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -34039.4
$ws.Range("H73").Value = 33499.4
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 33499.4
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 33499.4
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -35371.4
$ws.Range("H76").Value = 65000
$ws.Range("J76").Value = 65000
$ws.Range("L76").Value = 65000
$ws.Range("N76").Value = -65676
$ws.Range("H79").Value = 65000
$ws.Range("J79").Value = 65000
$ws.Range("L79").Value = 65000
$ws.Range("N79").Value = -67340
$ws.Range("H93").Value = 1148.6757
$ws.Range("I93").Value = 1099.12
$ws.Range("J93").Value = 1251.9166
$ws.Range("K93").Value = 1099.12
$ws.Range("L93").Value = 1251.9166
$ws.Range("M93").Value = 148.8800000000001
$ws.Range("N93").Value = -3747.9166
$ws.Range("H113").Value = 1829.0344
$ws.Range("I113").Value = 1779.3334
$ws.Range("K113").Value = 1779.3334
$ws.Range("M113").Value = 390.6666
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 7817.8335
$ws.Range("I23").Value = 510
$ws.Range("J23").Value = 11471.75
$ws.Range("K23").Value = 510
$ws.Range("L23").Value = 11471.75
$ws.Range("M23").Value = -281
$ws.Range("N23").Value = -11929.75
$ws.Range("H112").Value = 29949.5
$ws.Range("J112").Value = 29949.5
$ws.Range("L112").Value = 29949.5
$ws.Range("N112").Value = -32903.5
